$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-9 (columns A:E)
# A = Actual, B = Predicted-ish raw, C = ..., D = Residual, E = Confidence
$data = @(
    @(2, 3733.67808,            3763.6118784,         29.9337984000008,   0.03232710018566534),
    @(1, 2446.20288,            2208.856826880001,    237.3460531199989,  0.004195580278799644),
    @(4, 2639.32416,            2966.825663999998,    327.5015039999976,  0.003044126093255291),
    @(7, 2768.07168,            3216.434918399999,    448.3632383999989,  0.002225371179806778),
    @(6, 5391.3024,             6426.110592000005,    1034.808192000004,  0.0009654297076654092),
    @(5, 4667.097600000001,     5749.381439999999,    1082.283839999998,  0.0009231190968380011),
    @(3, 14162.2272,            12403.21420800001,    1759.012991999991,  0.0005681776239979058),
    @(0, 4506.1632,             6281.269632000007,    1775.106432000007,  0.0005630293218824379)
)

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 5).Value = $r[4]
    $row++
}
